$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on changed Price/Volume cells to preserve exact
# textual formatting (these columns store formatted text, not real numbers).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = "26.766.87"
$ws.Range("D3").Value = "1.648.50"
$ws.Range("E3").Value = "  +1.45%  "
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").Value = "216.31"
$ws.Range("E5").Value = "  +1.75%  "
$ws.Range("E6").Value = "  +1.72%  "
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("E8").Value = "  +2.02%  "
$ws.Range("E9").Value = "  +0.88%  "
$ws.Range("D10").Value = "19.19"
$ws.Range("E10").Value = "  +2.50%  "
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").Value = "1.879.52"
$ws.Range("E12").Value = "  +1.51%  "
$ws.Range("D13").Value = "1.660.82"
$ws.Range("E13").Value = "  +1.54%  "
$ws.Range("E14").Value = "  +1.66%  "
$ws.Range("E15").Value = "  +2.22%  "
$ws.Range("D16").Value = "65.52"
$ws.Range("E16").Value = "  +1.13%  "
$ws.Range("D17").Value = "26.792.01"
$ws.Range("E17").Value = "  +1.23%  "
$ws.Range("D18").Value = "0.0₃0745"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("D19").Value = "218.86"
$ws.Range("E19").Value = "  +2.71%  "
$ws.Range("E21").Value = "  +1.94%  "
$ws.Range("E22").Value = "  +0.68%  "
$ws.Range("E23").Value = "  +16.51%  "
$ws.Range("D24").Value = "9.51"
$ws.Range("E24").Value = "  +2.66%  "
$ws.Range("D25").Value = "146.40"
$ws.Range("E25").Value = "  -1.39%  "
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("E27").Value = "  +0.70%  "
$ws.Range("D28").Value = "7.13"
$ws.Range("E28").Value = "  +4.45%  "
$ws.Range("E29").Value = "  +1.68%  "
$ws.Range("E30").Value = "  +2.02%  "
$ws.Range("E31").Value = "  +1.62%  "
$ws.Range("E32").Value = "  +1.21%  "
$ws.Range("D33").Value = "3.01"
$ws.Range("E33").Value = "  +2.57%  "
$ws.Range("D34").Value = "1.281.62"
$ws.Range("E34").Value = "  +4.96%  "
$ws.Range("D35").Value = "1.54"
$ws.Range("E35").Value = "  +3.72%  "
$ws.Range("E36").Value = "  +2.24%  "
$ws.Range("E37").Value = "  +3.72%  "
$ws.Range("D38").Value = "0.537"
$ws.Range("E38").Value = "  +6.54%  "
$ws.Range("D39").Value = "0.829"
$ws.Range("E39").Value = "  +4.54%  "
$ws.Range("E40").Value = "  +0.40%  "
$ws.Range("D41").Value = "0.814"
$ws.Range("E41").Value = "  +3.01%  "
$ws.Range("E42").Value = "  -0.73%  "
$ws.Range("E43").Value = "  +2.43%  "
$ws.Range("D44").Value = "1.789.95"
$ws.Range("E44").Value = "  +1.59%  "
$ws.Range("D45").Value = "92.01"
$ws.Range("E45").Value = "  -0.76%  "
$ws.Range("D46").Value = "59.82"
$ws.Range("E46").Value = "  +9.52%  "
$ws.Range("E47").Value = "  +2.33%  "
$ws.Range("E48").Value = "  +1.29%  "
$ws.Range("D49").Value = "7.79"
$ws.Range("E49").Value = "  +3.92%  "
$ws.Range("E50").Value = "  +2.34%  "
$ws.Range("D51").Value = "0.407"
$ws.Range("E51").Value = "  +0.27%  "
